# Running sims on the new ac3 cluster
# Updates a few simulation-start-date values and appends a new log row (Set_14).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing simulation start dates (serial date numbers, keep existing date formatting) ---
$ws.Range("D14").Value = 45677
$ws.Range("D15").Value = 45701
$ws.Range("D16").Value = 45714

# --- Append a new row (row 17) for Set_14 ---
$ws.Range("A17").Value = "Set_14"
$ws.Range("B17").Value = "With burnin (scale = 0.033, ml = 0.5, ml_expt = 2, mut_ratio = 0.02)"
$ws.Range("C17").Value = 100

# Match the date number formatting used by the other rows (numFmtId 15) before setting the value
$ws.Range("D17").NumberFormat = $ws.Range("D16").NumberFormat
$ws.Range("D17").Value = 45734

# --- Update the sheet view so the newly added row is in view/selected ---
[void]$ws.Range("D17").Select()
